$wb = $excel.ActiveWorkbook

# --- "Request" sheet: TestCaseNumber renamed to the "ProductRules" naming
#     scheme, plus a few Description tweaks for the 1975 / Core / Lite /
#     1924 rows. ---
$reqWs = $wb.Worksheets.Item("Request")

$reqWs.Range("A2").Value = "TC _001_ProductRules_60100000"
$reqWs.Range("A3").Value = "TC _002_ProductRules_60200000"
$reqWs.Range("A4").Value = "TC _003_ProductRules_1975"
$reqWs.Range("C4").Value = "Getting product rule for  product 1975"
$reqWs.Range("A5").Value = "TC _004_ProductRules_1960E"
$reqWs.Range("A6").Value = "TC _005_ProductRules_1960C"
$reqWs.Range("C6").Value = "Getting product rule for Core V6.0"
$reqWs.Range("A7").Value = "TC _006_ProductRules_1960L"
$reqWs.Range("C7").Value = "Getting product rule for Lite V6.0"
$reqWs.Range("A8").Value = "TC _007_ProductRules_1977"
$reqWs.Range("A9").Value = "TC  _008_ProductRules_1924"
$reqWs.Range("C9").Value = "Getting product rule for Elite V6.0(Product 1924)"

# --- "Response" sheet: same TestCaseNumber renaming in column A. ---
$respWs = $wb.Worksheets.Item("Response")

$respWs.Range("A2").Value = "TC _001_ProductRules_60100000"
$respWs.Range("A3").Value = "TC _002_ProductRules_60200000"
$respWs.Range("A4").Value = "TC _003_ProductRules_1975"
$respWs.Range("A5").Value = "TC _004_ProductRules_1960E"
$respWs.Range("A6").Value = "TC _005_ProductRules_1960C"
$respWs.Range("A7").Value = "TC _006_ProductRules_1960L"
$respWs.Range("A8").Value = "TC _007_ProductRules_1977"
$respWs.Range("A9").Value = "TC  _008_ProductRules_1924"
